# Updated cryptos list on Tue Jan  9 21:53:57 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with the latest scrape, and re-syncs three rows (48/49/51) whose
# coin/link moved position in the source feed (RocketPoolETH, ordi, THORChain).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: latest price & volume
$ws.Range("D2").Value = "45.906.82"
$ws.Range("E2").Value = "  -2.58%  "

# Row 3: latest price & volume
$ws.Range("D3").Value = "2.355.19"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4: latest price & volume
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5: latest price & volume
$ws.Range("D5").Value = "'299.88"
$ws.Range("E5").Value = "  -2.03%  "

# Row 6: latest price & volume
$ws.Range("D6").Value = "'98.08"
$ws.Range("E6").Value = "  +0.76%  "

# Row 7: latest price & volume
$ws.Range("D7").Value = "'0.570"

# Row 8: latest price & volume
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.18%  "

# Row 9: latest price & volume
$ws.Range("D9").Value = "'0.510"
$ws.Range("E9").Value = "  -5.02%  "

# Row 10: latest price & volume
$ws.Range("D10").Value = "'34.58"
$ws.Range("E10").Value = "  -3.45%  "

# Row 11: latest price & volume
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -3.20%  "

# Row 12: latest price & volume
$ws.Range("E12").Value = "  -4.82%  "

# Row 13: latest price & volume
$ws.Range("E13").Value = "  -1.59%  "

# Row 14: latest price & volume
$ws.Range("D14").Value = "2.713.08"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15: latest price & volume
$ws.Range("D15").Value = "2.358.86"
$ws.Range("E15").Value = "  +0.77%  "

# Row 16: latest price & volume
$ws.Range("D16").Value = "'13.69"
$ws.Range("E16").Value = "  -3.44%  "

# Row 17: latest price & volume
$ws.Range("D17").Value = "'0.805"
$ws.Range("E17").Value = "  -3.36%  "

# Row 18: latest price & volume
$ws.Range("D18").Value = "46.047.24"
$ws.Range("E18").Value = "  -1.91%  "

# Row 19: latest price & volume
$ws.Range("D19").Value = "'12.54"
$ws.Range("E19").Value = "  -7.50%  "

# Row 20: latest price & volume
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  +2.24%  "

# Row 21: latest price & volume
$ws.Range("D21").Value = "'5.91"
$ws.Range("E21").Value = "  -4.39%  "

# Row 22: latest price & volume
$ws.Range("D22").Value = "'65.76"
$ws.Range("E22").Value = "  -2.60%  "

# Row 23: latest price & volume
$ws.Range("D23").Value = "'244.00"
$ws.Range("E23").Value = "  -2.32%  "

# Row 24: latest price & volume
$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -6.00%  "

# Row 25: latest price & volume
$ws.Range("E25").Value = "  +0.08%  "

# Row 26: latest price & volume
$ws.Range("D26").Value = "'1.87"
$ws.Range("E26").Value = "  -5.71%  "

# Row 27: latest price & volume
$ws.Range("D27").Value = "'40.47"
$ws.Range("E27").Value = "  -5.73%  "

# Row 28: latest price & volume
$ws.Range("E28").Value = "  -1.77%  "

# Row 29: latest price & volume
$ws.Range("D29").Value = "'9.63"
$ws.Range("E29").Value = "  -2.62%  "

# Row 30: latest price & volume
$ws.Range("D30").Value = "'20.51"
$ws.Range("E30").Value = "  +1.45%  "

# Row 31: latest price & volume
$ws.Range("D31").Value = "'3.54"
$ws.Range("E31").Value = "  +11.39%  "

# Row 32: latest price & volume
$ws.Range("D32").Value = "'2.80"
$ws.Range("E32").Value = "  +6.63%  "

# Row 33: latest price & volume
$ws.Range("D33").Value = "'144.78"
$ws.Range("E33").Value = "  -2.24%  "

# Row 34: latest price & volume
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  -7.61%  "

# Row 35: latest price & volume
$ws.Range("D35").Value = "'0.0770"
$ws.Range("E35").Value = "  -5.73%  "

# Row 36: latest price & volume
$ws.Range("E36").Value = "  -2.80%  "

# Row 37: latest price & volume
$ws.Range("D37").Value = "'0.115"
$ws.Range("E37").Value = "  -3.17%  "

# Row 38: latest price & volume
$ws.Range("D38").Value = "'1.76"
$ws.Range("E38").Value = "  -3.19%  "

# Row 39: latest price & volume
$ws.Range("D39").Value = "'15.26"
$ws.Range("E39").Value = "  +9.34%  "

# Row 40: latest price & volume
$ws.Range("D40").Value = "'3.85"
$ws.Range("E40").Value = "  -4.06%  "

# Row 41: latest price & volume
$ws.Range("E41").Value = "  -5.72%  "

# Row 42: latest price & volume
$ws.Range("D42").Value = "'3.12"
$ws.Range("E42").Value = "  -8.28%  "

# Row 43: latest price & volume
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.07%  "

# Row 44: latest price & volume
$ws.Range("D44").Value = "1.860.82"
$ws.Range("E44").Value = "  +1.60%  "

# Row 45: latest price & volume
$ws.Range("D45").Value = "'91.87"
$ws.Range("E45").Value = "  +3.48%  "

# Row 46: latest price & volume
$ws.Range("D46").Value = "'1.83"
$ws.Range("E46").Value = "  -7.49%  "

# Row 47: latest price & volume
$ws.Range("D47").Value = "'0.183"
$ws.Range("E47").Value = "  -6.20%  "

# Row 48: coin/link re-sync plus latest price & volume
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.587.44"
$ws.Range("E48").Value = "  +0.35%  "

# Row 49: coin/link re-sync plus latest price & volume
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'69.13"
$ws.Range("E49").Value = "  -8.30%  "

# Row 50: latest price & volume
$ws.Range("D50").Value = "'95.79"
$ws.Range("E50").Value = "  -3.26%  "

# Row 51: coin/link re-sync plus latest price & volume
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'4.72"
$ws.Range("E51").Value = "  -3.06%  "
